$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# N2: refreshed OAuth token payload (access_token / refresh_token / expires_in updated)
$ws.Range("N2").Value = "{'access_token': 'eeef7aa9-4555-4f1a-be0e-aecb559e1ded', 'token_type': 'bearer', 'refresh_token': '11aeea9c-abdb-4ee5-89c8-35ec6cba9f65', 'expires_in': 31551, 'scope': 'server', 'tenant_id': 'TAX', 'license': 'made by ling', 'eName': 'edmspicMic', 'staffNo': '1002', 'user_id': '9e2f885f-fccc-4f6c-9df7-5f42b421b7db', 'cName': '自动化用户002', 'active': True, 'dept_id': None, 'username': 'autoMic'}"

# L3: test result flips from PASS to FAIL, needs a new bold/red "宋体" font style
$ws.Range("L3").Value = "FAIL"
$ws.Range("L3").Font.Name = "宋体"
$ws.Range("L3").Font.Bold = $true
$ws.Range("L3").Font.Color = 255

# N3: response body now reflects the token failure instead of the big project-query payload
$ws.Range("N3").Value = "{'code': 40025003, 'message': 'Token不存在或验证错误', 'data': '8b0dc6ca-71d1-4581-9942-4c30b16ce49d', 'success': False, 'msg': 'Token不存在或验证错误'}"
